# Updates the "cryptos" price/volume sheet with refreshed values from the
# latest GitHub Actions run. Most rows just get new Price (D) / Volume(1h)
# (E) figures, but a handful of adjacent coin pairs swapped rank order
# (their whole row - Coin, Link, Price, Volume(1h) - is exchanged):
#   26/27 Cosmos<->Monero, 35/36 THORChain<->RenderToken,
#   41/42 VeChain<->Aave, 47/48 Cronos<->HuobiToken.
#
# Price values that look like plain numbers (e.g. "226.94") are forced to
# text via NumberFormat "@" before assignment so Excel keeps them as the
# original text strings instead of silently converting them to numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '39.736.67'
$ws.Range("E2").Value = '  +0.42%  '

$ws.Range("D3").Value = '2.176.35'
$ws.Range("E3").Value = '  +0.41%  '

$ws.Range("E4").Value = '  +0.06%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '226.94'
$ws.Range("E5").Value = '  -1.16%  '

$ws.Range("E6").Value = '  +0.48%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '63.15'
$ws.Range("E7").Value = '  -0.13%  '

$ws.Range("E8").Value = '  +0.08%  '

$ws.Range("E9").Value = '  -0.90%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0853'
$ws.Range("E10").Value = '  -0.90%  '

$ws.Range("E11").Value = '  +0.23%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '15.93'
$ws.Range("E12").Value = '  -1.16%  '

$ws.Range("D13").Value = '2.498.68'
$ws.Range("E13").Value = '  +0.72%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '21.85'
$ws.Range("E14").Value = '  -2.19%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.813'
$ws.Range("E15").Value = '  -1.22%  '

$ws.Range("E16").Value = '  -1.17%  '

$ws.Range("D17").Value = '2.176.90'
$ws.Range("E17").Value = '  +0.97%  '

$ws.Range("D18").Value = '39.708.14'

$ws.Range("D19").Value = '0.0₃0922'
$ws.Range("E19").Value = '  +7.83%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '71.81'
$ws.Range("E20").Value = '  -1.01%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.00'
$ws.Range("E21").Value = '  -2.67%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '229.79'
$ws.Range("E22").Value = '  +0.36%  '

$ws.Range("E23").Value = '  +0.01%  '

$ws.Range("E24").Value = '  -0.91%  '

$ws.Range("E25").Value = '  -0.80%  '

$ws.Range("B26").Value = 'Monero'
$ws.Range("C26").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '170.96'
$ws.Range("E26").Value = '  -1.18%  '

$ws.Range("B27").Value = 'Cosmos'
$ws.Range("C27").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.50'
$ws.Range("E27").Value = '  -2.02%  '

$ws.Range("E28").Value = '  +0.58%  '

$ws.Range("E29").Value = '  +1.24%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '19.86'
$ws.Range("E30").Value = '  +0.72%  '

$ws.Range("E31").Value = '  +3.74%  '

$ws.Range("E32").Value = '  -0.44%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.54'
$ws.Range("E33").Value = '  -2.68%  '

$ws.Range("E34").Value = '  -2.80%  '

$ws.Range("B35").Value = 'RenderToken'
$ws.Range("C35").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '3.93'
$ws.Range("E35").Value = '  +9.77%  '

$ws.Range("B36").Value = 'THORChain'
$ws.Range("C36").Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '6.97'
$ws.Range("E36").Value = '  -1.95%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.0618'
$ws.Range("E37").Value = '  -0.81%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.41'
$ws.Range("E38").Value = '  -2.06%  '

$ws.Range("E39").Value = '  +0.35%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '4.98'
$ws.Range("E40").Value = '  +16.91%  '

$ws.Range("B41").Value = 'Aave'
$ws.Range("C41").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '102.58'
$ws.Range("E41").Value = '  -1.33%  '

$ws.Range("B42").Value = 'VeChain'
$ws.Range("C42").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.0229'
$ws.Range("E42").Value = '  -1.42%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '17.75'
$ws.Range("E43").Value = '  -2.12%  '

$ws.Range("E44").Value = '  +2.60%  '

$ws.Range("D45").Value = '1.512.21'
$ws.Range("E45").Value = '  -1.33%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '7.90'
$ws.Range("E46").Value = '  +1.41%  '

$ws.Range("B47").Value = 'HuobiToken'
$ws.Range("C47").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.81'
$ws.Range("E47").Value = '  +0.14%  '

$ws.Range("B48").Value = 'Cronos'
$ws.Range("C48").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0921'
$ws.Range("E48").Value = '  -0.76%  '

$ws.Range("E49").Value = '  -1.59%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.000194'
$ws.Range("E50").Value = '  +30.76%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '49.67'
$ws.Range("E51").Value = '  +6.45%  '
